# Adding more diagnostic plots
# Clear the diagnostic columns (K:O -- LOOIC, n.big.Rhat, n.bad.Pareto,
# p.bad.Pareto, Nhats) for the model rows, leaving only the model name (J)
# in place, for every data row in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(3, 4, 5, 6, 10, 11, 12, 13, 15, 16, 17, 18, 22, 23, 24, 25, 27, 28, 29, 30, 34, 35, 36, 37)

foreach ($r in $rows) {
    $ws.Range("K" + $r + ":O" + $r).ClearContents()
}

$ws.Range("A38").Select() | Out-Null
